# SIMS year division script + quick fix for false year placements
# Adds a missing otolith record (G-BH-91588 / #G) as a new row 30, and
# updates the sheet's active selection/scroll to reflect where the user
# was working when they made the fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quick fix: add the missing "G" otolith entry that had been placed
# under the wrong year previously.
$ws.Range("A30").Value = "G-BH-91588"
$ws.Range("B30").Value = "#G"
$ws.Range("C30").Value = 3607
$ws.Range("D30").Value = 1200

# Reflect the place the user was last looking/selecting in the sheet.
$ws.Range("D29:D30").Select()
